# Weekly update: a new week's Albahaca price record is reported for
# "Agrícola del Norte S.A. de Arica". The new observation is inserted at
# the top of the data block (row 33), pushing all existing rows (33-75)
# down by one (to 34-76).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 33; rows 33..75 shift down to 34..76 and the sheet
# dimension grows from A1:R75 to A1:R76 automatically.
$ws.Rows.Item(33).Insert()

# Populate the newly inserted row 33 with this week's record. Most of the
# descriptive/price columns repeat the prior top record's values; only the
# date (D) and volume (J) are new for this observation.
$ws.Cells.Item(33, 1).Value  = 1
$ws.Cells.Item(33, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(33, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(33, 4).Value  = 45195
$ws.Cells.Item(33, 5).Value  = 15
$ws.Cells.Item(33, 6).Value  = 100112052
$ws.Cells.Item(33, 7).Value  = "Albahaca"
$ws.Cells.Item(33, 8).Value  = "Sin especificar"
$ws.Cells.Item(33, 9).Value  = "Primera"
$ws.Cells.Item(33, 10).Value = 300
$ws.Cells.Item(33, 11).Value = 900
$ws.Cells.Item(33, 12).Value = 1000
$ws.Cells.Item(33, 13).Value = 950
$ws.Cells.Item(33, 14).Value = "$/paquete"
$ws.Cells.Item(33, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(33, 16).Value = 950
$ws.Cells.Item(33, 17).Value = 1
$ws.Cells.Item(33, 18).Value = "Hortaliza"
